$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as text (matches original inlineStr cells)

$ws.Range("D2").Value = "'67.973.57"
$ws.Range("E2").Value = "'  -0.13%  "

$ws.Range("D3").Value = "'3.556.70"
$ws.Range("E3").Value = "'  -1.68%  "

$ws.Range("E4").Value = "'  +0.07%  "

$ws.Range("D5").Value = "'205.55"
$ws.Range("E5").Value = "'  +6.01%  "

$ws.Range("D6").Value = "'556.91"
$ws.Range("E6").Value = "'  -5.13%  "

$ws.Range("D7").Value = "'3.545.14"
$ws.Range("E7").Value = "'  -1.85%  "

$ws.Range("D8").Value = "'0.609"
$ws.Range("E8").Value = "'  -1.75%  "

$ws.Range("E9").Value = "'  +0.02%  "

$ws.Range("D10").Value = "'0.665"
$ws.Range("E10").Value = "'  -2.38%  "

$ws.Range("D11").Value = "'64.30"
$ws.Range("E11").Value = "'  +15.32%  "

$ws.Range("D12").Value = "'0.145"
$ws.Range("E12").Value = "'  -4.80%  "

$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = "'  -6.18%  "

$ws.Range("D14").Value = "'9.98"
$ws.Range("E14").Value = "'  -0.35%  "

$ws.Range("D15").Value = "'4.119.26"
$ws.Range("E15").Value = "'  -1.81%  "

$ws.Range("D16").Value = "'3.546.64"
$ws.Range("E16").Value = "'  -2.13%  "

$ws.Range("E17").Value = "'  -1.09%  "

$ws.Range("D18").Value = "'18.74"
$ws.Range("E18").Value = "'  +1.21%  "

$ws.Range("D19").Value = "'67.689.98"
$ws.Range("E19").Value = "'  -0.36%  "

$ws.Range("D20").Value = "'12.01"
$ws.Range("E20").Value = "'  -4.29%  "

$ws.Range("D21").Value = "'1.05"
$ws.Range("E21").Value = "'  -3.15%  "

$ws.Range("D22").Value = "'397.50"
$ws.Range("E22").Value = "'  -1.98%  "

$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "'  -4.60%  "

$ws.Range("D24").Value = "'12.13"
$ws.Range("E24").Value = "'  -10.62%  "

$ws.Range("D25").Value = "'83.35"
$ws.Range("E25").Value = "'  -3.19%  "

$ws.Range("D26").Value = "'2.84"
$ws.Range("E26").Value = "'  -3.91%  "

$ws.Range("D27").Value = "'12.29"
$ws.Range("E27").Value = "'  -2.91%  "

$ws.Range("D28").Value = "'3.79"
$ws.Range("E28").Value = "'  -3.49%  "

$ws.Range("D29").Value = "'9.02"
$ws.Range("E29").Value = "'  -2.04%  "

$ws.Range("D30").Value = "'31.18"
$ws.Range("E30").Value = "'  -1.11%  "

$ws.Range("B31").Value = "'Bittensor"
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'695.71"
$ws.Range("E31").Value = "'  +1.65%  "

$ws.Range("B32").Value = "'NEARProtocol"
$ws.Range("C32").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.35"
$ws.Range("E32").Value = "'  -11.85%  "

$ws.Range("D33").Value = "'11.94"
$ws.Range("E33").Value = "'  -2.63%  "

$ws.Range("D34").Value = "'64.44"
$ws.Range("E34").Value = "'  -0.23%  "

$ws.Range("D35").Value = "'0.112"
$ws.Range("E35").Value = "'  -5.08%  "

$ws.Range("D36").Value = "'39.76"
$ws.Range("E36").Value = "'  -6.66%  "

$ws.Range("D37").Value = "'0.416"
$ws.Range("E37").Value = "'  -1.58%  "

$ws.Range("E38").Value = "'  -0.09%  "

$ws.Range("D39").Value = "'3.07"
$ws.Range("E39").Value = "'  -1.92%  "

$ws.Range("D40").Value = "'0.133"
$ws.Range("E40").Value = "'  -1.00%  "

$ws.Range("D41").Value = "'3.105.72"
$ws.Range("E41").Value = "'  -2.81%  "

$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "'  -0.12%  "

$ws.Range("D43").Value = "'0.0₃0709"
$ws.Range("E43").Value = "'  -10.06%  "

$ws.Range("B44").Value = "'dogwifhat"
$ws.Range("C44").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.92"
$ws.Range("E44").Value = "'  +14.48%  "

$ws.Range("B45").Value = "'Fetch.AI"
$ws.Range("C45").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.60"
$ws.Range("E45").Value = "'  -12.99%  "

$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = "'  +7.53%  "

$ws.Range("D47").Value = "'0.0405"
$ws.Range("E47").Value = "'  -4.03%  "

$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "'  -2.66%  "

$ws.Range("B49").Value = "'THORChain"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'8.47"
$ws.Range("E49").Value = "'  -4.15%  "

$ws.Range("D50").Value = "'138.11"
$ws.Range("E50").Value = "'  -3.99%  "

$ws.Range("B51").Value = "'ApeXProtocol"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'2.98"
$ws.Range("E51").Value = "'  -3.55%  "
